$d = $word.ActiveDocument

$replacements = @(
    @("2024-03-05 Tuesday", "2024-03-06 Wednesday"),
    @("76×65=4940", "57×31=1767"),
    @("34×17=578", "55×80=4400"),
    @("41×12=492", "16×49=784"),
    @("95×96=9120", "58×32=1856"),
    @("48×40=1920", "48×76=3648"),
    @("85×80=6800", "73×18=1314"),
    @("60×22=1320", "11×31=341"),
    @("21×40=840", "83×40=3320"),
    @("44×36=1584", "63×73=4599"),
    @("82×80=6560", "69×40=2760"),
    @("74×43=3182", "66×73=4818"),
    @("52×14=728", "18×40=720"),
    @("20×18=360", "77×45=3465"),
    @("20×15=300", "63×42=2646"),
    @("81×87=7047", "40×46=1840"),
    @("18×74=1332", "19×84=1596"),
    @("59×72=4248", "39×64=2496"),
    @("39×35=1365", "46×22=1012"),
    @("70×80=5600", "77×14=1078"),
    @("81×28=2268", "24×98=2352"),
    @("84×86=7224", "83×75=6225"),
    @("75×57=4275", "70×69=4830"),
    @("60×88=5280", "18×87=1566"),
    @("74×50=3700", "62×34=2108"),
    @("53×68=3604", "65×51=3315")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
